$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 378, shifting existing rows 378:501 down to 379:502
$ws.Rows.Item(378).Insert()

# Populate the newly inserted row 378 with the new data record
$ws.Cells.Item(378, 1).Value = 6
$ws.Cells.Item(378, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(378, 3).Value = "Metropolitana"
$ws.Cells.Item(378, 4).Value = 44627
$ws.Cells.Item(378, 5).Value = 13
$ws.Cells.Item(378, 6).Value = 100112012
$ws.Cells.Item(378, 7).Value = "Espinaca"
$ws.Cells.Item(378, 8).Value = "Sin especificar"
$ws.Cells.Item(378, 9).Value = "Primera"
$ws.Cells.Item(378, 10).Value = 370
$ws.Cells.Item(378, 11).Value = 9000
$ws.Cells.Item(378, 12).Value = 10000
$ws.Cells.Item(378, 13).Value = 9405
$ws.Cells.Item(378, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(378, 15).Value = "Región Metropolitana"
$ws.Cells.Item(378, 16).Value = 940
$ws.Cells.Item(378, 17).Value = 10
$ws.Cells.Item(378, 18).Value = "Hortaliza"
